$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update M2: KODE_FIXED_INCOME value changes from OBL00107 to OBL00108
$ws.Range("M2").Value = "OBL00108"

# Update F2: the last line of the preparation note changes from
# "Kode Fixed Income : Hasil Generate" to "Kode Fixed Income : OBL00108"
$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 20/21 - Analis Investasi/Asisten Investasi;`nKode Fixed Income : OBL00108"

# Update the selected/active cell in the sheet view to G2 (was N2), and reset the top-left visible cell
$ws.Range("G2").Select()
